# Apply the "4 - Botões da home não funcionam / 9 - Validar senha do usuário"
# update to the bug tracker workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the formatting of an existing "visible" bug row (row 3) down onto the
# two new rows (17 and 18) so they pick up the same styles (s="6"/"7"/"8")
# used by the other non-hidden rows.
$ws.Range("A3:C3").Copy()
$ws.Range("A17:C18").PasteSpecial(-4122)

# Row 17: Bug #4 - Botões da Home não funcionam (assigned to Paulo / Tarcísio)
$ws.Cells.Item(17, 1).Value = 4
$ws.Cells.Item(17, 2).Value = "Botões da Home não funcionam"
$ws.Cells.Item(17, 3).Value = "Paulo / Tarcísio"

# Row 18: Bug #9 - Validar antiga senha do usuário ao editar perfil e trocar senha
$ws.Cells.Item(18, 1).Value = 9
$ws.Cells.Item(18, 2).Value = "Validar antiga senha do usuário ao editar perfil e trocar senha"
$ws.Cells.Item(18, 3).Value = "Paulo / Tarcísio"

# Give row 18 its slightly taller row height, matching the diff (ht="30").
$ws.Rows.Item(18).RowHeight = 30

# Grow the two tables (ListObjects) to include the two new rows.
$loBugs = $ws.ListObjects.Item(1)
$loBugs.Resize($ws.Range("B1:C18"))

$loNums = $ws.ListObjects.Item(2)
$loNums.Resize($ws.Range("A1:A18"))

# Update the active selection to match what was left selected in the file.
[void]$ws.Range("D21").Select()
